$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column -> new value mapping (applies identically to rows 2 and 3)
$updates = @{
    "D"  = 0.133
    "E"  = 0.224
    "K"  = 14.9
    "L"  = 0.5665399239543726
    "M"  = 7.84
    "N"  = 0.1058029689608637
    "O"  = 0.5261744966442953
    "P"  = 7.84
    "Q"  = 0.1058029689608637
    "R"  = 0.5261744966442953
    "U"  = 301.6
    "V"  = 4.070175438596491
    "W"  = 0.3260393873085339
    "X"  = 0.1115604283921338
    "Y"  = 0.2144789589164001
    "Z"  = 0.267005076142132
    "AB" = 0.07155442695350728
    "AC" = -0.07155442695350728
    "AD" = 152.4
    "AF" = 152.4
    "AG" = -149.2
    "AH" = 0.6728476821192053
    "AI" = 0.7437774524158126
    "AJ" = 1.986684420772303
    "AK" = 1.542916235780765
}

foreach ($row in 2, 3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}

$wb.Save()
